# Applies:
#   1. Refresh the cached "datetimeFigureOut" footer field (slide master +
#      every slide layout) from 12/3/2015 -> 12/5/2015.
#   2. Fix the "Idisposable" -> "IDisposable" typo in the code sample.

$p = $ppt.ActivePresentation

$oldDate = "12/3/2015"
$newDate = "12/5/2015"
$ppPlaceholderDate = 16

function Update-DatePlaceholders {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) {
            continue
        }
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if (-not $isDatePlaceholder) {
            continue
        }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# The date footer placeholder lives on the slide master and on every
# slide layout that inherits from it.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# Fix the "Idisposable" typo wherever it shows up in the deck, while
# touching only the misspelled substring so the surrounding run
# formatting (font, color, etc.) is left untouched.
$typo = "Idisposable"
$fix = "IDisposable"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) {
            continue
        }
        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text
        $idx = $text.IndexOf($typo)
        while ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $typo.Length)
            $sub.Text = $fix
            $text = $tr.Text
            $idx = $text.IndexOf($typo)
        }
    }
}
